{"js": "// Replace the merge-field reference \"hearingLocation.venue_name\" with\n// \"caseManagementLocation.venue_name\" inside the Heading1 paragraph that\n// reads: <<cs_{writtenByJudge}>><<hearingLocation.venue_name>><<else>> Online Civil Claims<<es_>>\nconst body = context.document.body;\n\nconst results = body.search(\"hearingLocation.venue_name\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"caseManagementLocation.venue_name\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the merge-field reference \"hearingLocation.venue_name\" with\n# \"caseManagementLocation.venue_name\" inside the Heading1 paragraph that\n# reads: <<cs_{writtenByJudge}>><<hearingLocation.venue_name>><<else>> Online Civil Claims<<es_>>\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"hearingLocation.venue_name\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"caseManagementLocation.venue_name\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $find.MatchSoundsLike, $find.MatchAllWordForms, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
